# Auto-generated script to apply scheduled runner updates to market-price/profit data
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I15").Value = 2595.8235
$ws.Range("K15").Value = 7787.470499999999
$ws.Range("M15").Value = -7618.470499999999
$ws.Range("H15").Value = 2595.8235
$ws.Range("L17").Value = 7592.400000000001
$ws.Range("J17").Value = 2530.8
$ws.Range("H17").Value = 2530.8
$ws.Range("N17").Value = -7928.400000000001
$ws.Range("I40").Value = 7525.857
$ws.Range("K40").Value = 7525.857
$ws.Range("M40").Value = -7350.857
$ws.Range("H40").Value = 6997.0835
$ws.Range("J41").Value = 3246.5
$ws.Range("M41").Value = 202.9375
$ws.Range("I41").Value = 237.0625
$ws.Range("K41").Value = 237.0625
$ws.Range("L41").Value = 3246.5
$ws.Range("H41").Value = 571.44446
$ws.Range("N41").Value = -4126.5
$ws.Range("L43").Value = 1770
$ws.Range("H43").Value = 4405.857
$ws.Range("N43").Value = -1908
$ws.Range("J43").Value = 1770
$ws.Range("I86").Value = 142858110
$ws.Range("K86").Value = 142858110
$ws.Range("M86").Value = -142856987
$ws.Range("H86").Value = 71430470
$ws.Range("I89").Value = 142858110
$ws.Range("K89").Value = 714290550
$ws.Range("M89").Value = -714284934
$ws.Range("H89").Value = 71430470
$ws.Range("L112").Value = 6773.625
$ws.Range("J112").Value = 2257.875
$ws.Range("H112").Value = 2257.875
$ws.Range("N112").Value = -8989.625
$ws.Range("N113").Value = -12070.5
$ws.Range("I113").Value = 5252
$ws.Range("K113").Value = 5252
$ws.Range("L113").Value = 5562.5
$ws.Range("J113").Value = 5562.5
$ws.Range("M113").Value = -1998
$ws.Range("H113").Value = 5459
$ws.Range("M116").Value = -994492.9
$ws.Range("H116").Value = 873617.75
$ws.Range("I116").Value = 997934.9
$ws.Range("K116").Value = 997934.9
$ws.Range("K131").Value = 5013586.800000001
$ws.Range("L131").Value = 10200
$ws.Range("M131").Value = -5008546.800000001
$ws.Range("H131").Value = 1115263.8
$ws.Range("N131").Value = -20280
$ws.Range("J131").Value = 3400
$ws.Range("I131").Value = 1671195.6
$ws.Range("I135").Value = 323.53845
$ws.Range("K135").Value = 2911.84605
$ws.Range("M135").Value = -376.8460500000001
$ws.Range("H135").Value = 670.8
$ws.Range("J137").Value = 4970.852
$ws.Range("L137").Value = 14912.556
$ws.Range("H137").Value = 2993.75
$ws.Range("N137").Value = -20012.556
$ws.Range("L138").Value = 1282231.02
$ws.Range("H138").Value = 281402.72
$ws.Range("J138").Value = 427410.34
$ws.Range("M138").Value = -5411.6774
$ws.Range("N138").Value = -1292511.02
$ws.Range("I138").Value = 3517.2258
$ws.Range("K138").Value = 10551.6774
$ws.Range("L141").Value = 14698.9995
$ws.Range("H141").Value = 4874.625
$ws.Range("J141").Value = 4899.6665
$ws.Range("M141").Value = -9218.5
$ws.Range("N141").Value = -25058.9995
$ws.Range("I141").Value = 4799.5
$ws.Range("K141").Value = 14398.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -3952.86
$ws.Range("H32").Value = 5159.75
$ws.Range("N32").Value = -13399.5
$ws.Range("J32").Value = 12825.5
$ws.Range("I32").Value = 4239.86
$ws.Range("K32").Value = 4239.86
$ws.Range("L32").Value = 12825.5
$ws.Range("M74").Value = -370983.88
$ws.Range("H74").Value = 204159.64
$ws.Range("I74").Value = 371857.88
$ws.Range("K74").Value = 371857.88
$ws.Range("I77").Value = 371857.88
$ws.Range("K77").Value = 1859289.4
$ws.Range("M77").Value = -1854921.4
$ws.Range("H77").Value = 204159.64
$ws.Range("J102").Value = 4721.6665
$ws.Range("M102").Value = -3332
$ws.Range("N102").Value = -7965.6665
$ws.Range("L102").Value = 4721.6665
$ws.Range("I102").Value = 4954
$ws.Range("K102").Value = 4954
$ws.Range("H102").Value = 4920.8096
$ws.Range("H132").Value = 2278.4167
$ws.Range("I132").Value = 1685
$ws.Range("K132").Value = 5055
$ws.Range("M132").Value = -2525

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M20").Value = -15628804
$ws.Range("H20").Value = 13161790
$ws.Range("I20").Value = 15629051
$ws.Range("K20").Value = 15629051
$ws.Range("H111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 3450.44
$ws.Range("K31").Value = 3450.44
$ws.Range("M31").Value = -3155.44
$ws.Range("H31").Value = 4137.162
$ws.Range("I34").Value = 3450.44
$ws.Range("K34").Value = 3450.44
$ws.Range("M34").Value = -3248.44
$ws.Range("H34").Value = 4137.162
$ws.Range("L124").Value = 42326
$ws.Range("H124").Value = 42326
$ws.Range("N124").Value = -47236
$ws.Range("J124").Value = 42326
$ws.Range("H132").Value = 11366378
$ws.Range("N132").Value = -20059.1432
$ws.Range("J132").Value = 4999.7144
$ws.Range("I132").Value = 13515827
$ws.Range("K132").Value = 40547481
$ws.Range("L132").Value = 14999.1432
$ws.Range("M132").Value = -40544951
$ws.Range("I134").Value = 3392
$ws.Range("K134").Value = 10176
$ws.Range("M134").Value = -7641
$ws.Range("H134").Value = 3956.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L82").Value = 34470
$ws.Range("M82").Value = -20636
$ws.Range("H82").Value = 10744
$ws.Range("N82").Value = -35282
$ws.Range("J82").Value = 11490
$ws.Range("I82").Value = 7014
$ws.Range("K82").Value = 21042
$ws.Range("I85").Value = 7014
$ws.Range("K85").Value = 21042
$ws.Range("L85").Value = 34470
$ws.Range("M85").Value = -19638
$ws.Range("H85").Value = 10744
$ws.Range("N85").Value = -37278
$ws.Range("J85").Value = 11490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M97").Value = -326.5333000000001
$ws.Range("H97").Value = 2323.4
$ws.Range("N97").Value = -5566.7
$ws.Range("J97").Value = 4574.7
$ws.Range("I97").Value = 822.5333000000001
$ws.Range("K97").Value = 822.5333000000001
$ws.Range("L97").Value = 4574.7
$ws.Range("J102").Value = 21681.4
$ws.Range("M102").Value = 813.51166
$ws.Range("N102").Value = -24925.4
$ws.Range("L102").Value = 21681.4
$ws.Range("I102").Value = 808.48834
$ws.Range("K102").Value = 808.48834
$ws.Range("H102").Value = 4746.7734
$ws.Range("L112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("I113").Value = 4117.654
$ws.Range("K113").Value = 4117.654
$ws.Range("M113").Value = -1947.654
$ws.Range("H113").Value = 4157.1562
$ws.Range("H132").Value = 2417
$ws.Range("N132").Value = -11555
$ws.Range("J132").Value = 2165
$ws.Range("I132").Value = 2495.2068
$ws.Range("K132").Value = 7485.6204
$ws.Range("L132").Value = 6495
$ws.Range("M132").Value = -4955.6204
$ws.Range("L136").Value = 28100.571
$ws.Range("H136").Value = 9366.857
$ws.Range("N136").Value = -33200.571
$ws.Range("J136").Value = 9366.857
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K68").Value = 7999.75
$ws.Range("L68").Value = 8501.5
$ws.Range("M68").Value = -7250.75
$ws.Range("H68").Value = 8167
$ws.Range("N68").Value = -9999.5
$ws.Range("J68").Value = 8501.5
$ws.Range("I68").Value = 7999.75
$ws.Range("M71").Value = -36254.75
$ws.Range("H71").Value = 8167
$ws.Range("N71").Value = -49995.5
$ws.Range("J71").Value = 8501.5
$ws.Range("I71").Value = 7999.75
$ws.Range("K71").Value = 39998.75
$ws.Range("L71").Value = 42507.5
$ws.Range("H132").Value = 2889.9883
$ws.Range("I132").Value = 2183.6973
$ws.Range("K132").Value = 6551.091899999999
$ws.Range("M132").Value = -4021.091899999999
$ws.Range("I136").Value = 4161.8423
$ws.Range("K136").Value = 12485.5269
$ws.Range("M136").Value = -9935.526900000001
$ws.Range("H136").Value = 4344.276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 62500996
$ws.Range("J100").Value = 142858290
$ws.Range("M100").Value = -1231.4444
$ws.Range("N100").Value = -285717662
$ws.Range("I100").Value = 886.2222
$ws.Range("K100").Value = 1772.4444
$ws.Range("L100").Value = 285716580
$ws.Range("L107").Value = 1112.1429
$ws.Range("M107").Value = -278.5712999999996
$ws.Range("H107").Value = 612.1429000000001
$ws.Range("N107").Value = -4952.1429
$ws.Range("J107").Value = 370.7143
$ws.Range("I107").Value = 732.8570999999999
$ws.Range("K107").Value = 2198.5713
$ws.Range("M122").Value = -2472.2173
$ws.Range("H122").Value = 8335389.5
$ws.Range("I122").Value = 1640.7391
$ws.Range("K122").Value = 4922.2173
$ws.Range("H132").Value = 8549843
$ws.Range("N132").Value = -10816.1426
$ws.Range("J132").Value = 1918.7142
$ws.Range("I132").Value = 10419701
$ws.Range("K132").Value = 31259103
$ws.Range("L132").Value = 5756.142599999999
$ws.Range("M132").Value = -31256573
$ws.Range("N135").Value = -64398.168
$ws.Range("J135").Value = 54258.168
$ws.Range("L135").Value = 54258.168
$ws.Range("H135").Value = 54258.168
$ws.Range("L136").Value = 21281.6661
$ws.Range("H136").Value = 18870002
$ws.Range("N136").Value = -26381.6661
$ws.Range("J136").Value = 7093.8887

